$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update timestamp (row 1)
$ws.Range("A1").Value = "Datos actualizados a 15 de Octubre de 2020 a las 06:27"

# Row 5: India - updated case numbers
$ws.Range("B5").Value = 7307097
$ws.Range("C5").Value = 2027
$ws.Range("D5").Value = 6383441
$ws.Range("E5").Value = 812345

# Rows 54-56: Honduras moves up above Bielorrusia & Venezuela (with updated data),
# Bielorrusia and Venezuela shift down one row each (values unchanged).
$ws.Range("A54").Value = "Honduras"
$ws.Range("B54").Value = 85458
$ws.Range("C54").Value = 606
$ws.Range("D54").Value = 32990
$ws.Range("E54").Value = 49935
$ws.Range("F54").Value = 0
$ws.Range("G54").Value = 5
$ws.Range("H54").Value = 2533

$ws.Range("A55").Value = "Bielorrusia"
$ws.Range("B55").Value = 85121
$ws.Range("C55").Value = 0
$ws.Range("D55").Value = 78218
$ws.Range("E55").Value = 5992
$ws.Range("F55").Value = 0
$ws.Range("G55").Value = 0
$ws.Range("H55").Value = 911

$ws.Range("A56").Value = "Venezuela"
$ws.Range("B56").Value = 85005
$ws.Range("C56").Value = 0
$ws.Range("D56").Value = 76262
$ws.Range("E56").Value = 8029
$ws.Range("F56").Value = 0
$ws.Range("G56").Value = 0
$ws.Range("H56").Value = 714

# Row 144: Tailandia - updated case numbers
$ws.Range("B144").Value = 3665
$ws.Range("C144").Value = 13
$ws.Range("D144").Value = 3463
$ws.Range("E144").Value = 143

# Row 187: Butan - updated case numbers
$ws.Range("B187").Value = 316
$ws.Range("C187").Value = 3
$ws.Range("D187").Value = 294
$ws.Range("E187").Value = 22
